$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.113.70"
$ws.Range("E2").Value = "  -2.06%  "

$ws.Range("D3").Value = "1.865.73"
$ws.Range("E3").Value = "  -1.96%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.00%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5135"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3761"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07135"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8886"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.38%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07546"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.46%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.860.95"
$ws.Range("E13").Value = "  -2.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.308"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.56%  "

$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008470"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.61%  "

$ws.Range("E18").Value = "  -2.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "27.153.26"
$ws.Range("E20").Value = "  -2.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.71%  "

$ws.Range("D22").Value = "2.091.79"
$ws.Range("E22").Value = "  -1.58%  "

$ws.Range("E23").Value = "  -3.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.445"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.837"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.087"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.64%  "

$ws.Range("E30").Value = "  -4.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.666"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09192"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05102"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.077"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.155"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7226"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02035"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.40%  "

$ws.Range("E38").Value = "  +0.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.492"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.075"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5271"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.487"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1466"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4615"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.952"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.560"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.97%  "
